# BUCK.mac slide: split the "802" material-number literal and the
# "!General Aluminum Beam, A = 1.1, " comment run into separate runs so
# the macro argument ("44") and the trailing ", A = 1.1, " text can be
# combined / edited independently later ("ability to combine multiple
# runs for BUCK.mac").
#
#   ", 10, 1.1, 4.25, 8.1, 8.1, 802   "          -> ", 10, 1.1, 4.25, 8.1, 8.1, " + "44 "
#   "!General Aluminum Beam, A = 1.1, "          -> "!Material 44 Beam" + ", A = 1.1, "

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(13)
$shape = $s.Shapes.Item(13)
$tr = $shape.TextFrame.TextRange

# Locate the two runs of interest by searching for their current text so
# the script doesn't depend on hard-coded character offsets.
$fullText = $tr.Text
$oldRun2Text = ", 10, 1.1, 4.25, 8.1, 8.1, 802   "
$oldRun3Text = "!General Aluminum Beam, A = 1.1, "

$run2Start = $fullText.IndexOf($oldRun2Text) + 1   # TextRange positions are 1-based
$run3Start = $fullText.IndexOf($oldRun3Text) + 1

# --- Run 2 (white "bg1" text): ", 10, 1.1, 4.25, 8.1, 8.1, 802   " ---
# becomes ", 10, 1.1, 4.25, 8.1, 8.1, " + a new run "44 ".
$newRun2Text = ", 10, 1.1, 4.25, 8.1, 8.1, "
$matNumText = "44 "

$r2 = $tr.Characters($run2Start, $oldRun2Text.Length)
$r2.Text = $newRun2Text

$r2fresh = $tr.Characters($run2Start, $newRun2Text.Length)
$r2fresh.InsertAfter($matNumText) | Out-Null

# Touching the Font property on just the inserted text forces it into
# its own run (matching formatting, but a distinct <a:r>).
$matNumStart = $run2Start + $newRun2Text.Length
$matNumRange = $tr.Characters($matNumStart, $matNumText.Length)
$matNumRange.Font.Bold = $true

# --- Run 3 (green "00B050" text): "!General Aluminum Beam, A = 1.1, " ---
# becomes "!Material 44 Beam" + a new run ", A = 1.1, ".
$newRun3Text = "!Material 44 Beam"
$suffixText = ", A = 1.1, "

# run3's start shifts by however much run2's total length changed.
$run3StartNew = $run3Start + ($newRun2Text.Length + $matNumText.Length) - $oldRun2Text.Length

$r3 = $tr.Characters($run3StartNew, $oldRun3Text.Length)
$r3.Text = $newRun3Text

$r3fresh = $tr.Characters($run3StartNew, $newRun3Text.Length)
$r3fresh.InsertAfter($suffixText) | Out-Null

$suffixStart = $run3StartNew + $newRun3Text.Length
$suffixRange = $tr.Characters($suffixStart, $suffixText.Length)
$suffixRange.Font.Bold = $true
